# Load the improved microcode instruction names into the instruction set map.
# Rotate-left/right instructions renamed (RLN/RRN -> RLD/RRD) and the
# branch-if-zero opcode renamed to branch-if-equal-zero (BZ i16 -> BEZ i16).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F3").Value = "RLD"
$ws.Range("F5").Value = "RRD"
$ws.Range("M4").Value = "BEZ i16"

$ws.Range("M11").Select()
